# Final Project Narratives.docx edit script
# Implements:
#  1. Paragraph "Our example serializations..." (XML/JSON Serializations section):
#     - "One XML file is used for each of the top classes, " ->
#       "We chose to split the serialized data into three files, one for each of the top classes, "
#     - "crIssues" -> "crIssue" (first, italic, classes list)
#     - append a large new explanatory block with several italic terms.
#  2. Paragraph "For the XML serializations, ..." (External vocabularies section):
#     - "we chose to incorporate two standard external vocabularies, Dublin Core and Friend of a
#       Friend. First we set an imaginary default namespace in each XML document
#       (\u201chttps://example.com/maryse_and_joe_project/\u201d) with which to define all the
#       elements that aren\u2019t given namespaces. " ->
#       "we chose to incorporate bits and pieces of four authoritative external vocabularies. "
#  3. Header: normalize "Joseph Muller & Maryse " / "Lundering-Timpano" into a single run (no
#     visible text change).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "XML / JSON Serializations" paragraph
# ---------------------------------------------------------------------------

$p11 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Our example serializations are provided in XML*") {
        $p11 = $cand
        break
    }
}

# 1a. Expand "One XML file is used for each of the top classes, " into the new lead-in text.
$r = $p11.Range
$r.Find.Execute(
    "One XML file is used for each of the top classes, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We chose to split the serialized data into three files, one for each of the top classes, ",
    2)

# 1b. "crIssues" (the one right before ", congressPerson") -> "crIssue" (keeps italic formatting
#     of the run being replaced).
$r = $p11.Range
$r.Find.Execute("crIssues", $true, $false, $false, $false, $false, $true, 1, $false, "crIssue", 2)

# 1c. Append the new trailing content as plain (non-italic) text first; italics are applied
#     afterwards by searching back through the paragraph, so that newly inserted text never
#     inherits formatting from a preceding italic run.
$pEnd = $p11.Range.End
$ins = $d.Range($pEnd - 1, $pEnd - 1)
$ins.InsertAfter("In each file is a wrapper as the root element. For the RELAX NG validation schemas, we specify that one or more instances of the main element must exist. However, not every crIssue needs to have a document (speech) instance, since there might not be any speeches in some issues of the Congressional Record. We also chose to validate data types, as we have four different XSD data types in our data: string, anyURI, float, and integer. One RELAX NG file is provided for each XML file.")

# 1d. Apply italics to the relevant terms within the newly appended text. We restrict the search
#     to the paragraph range and walk forward sequentially to keep the correct term instances.
$r = $p11.Range
$r.Find.Execute("However, not every ")
$r.Collapse(0)
$r.Find.Execute("crIssue")
$r.Italic = 1

$r = $p11.Range
$r.Find.Execute(" needs to have a ")
$r.Collapse(0)
$r.Find.Execute("document")
$r.Italic = 1

$r = $p11.Range
$r.Find.Execute("some issues of the ")
$r.Collapse(0)
$r.Find.Execute("Congressional Record")
$r.Italic = 1

$r = $p11.Range
$r.Find.Execute("our data: ")
$r.Collapse(0)
$r.Find.Execute("string")
$r.Italic = 1
$r.Collapse(0)
$r.Find.Execute("anyURI")
$r.Italic = 1
$r.Collapse(0)
$r.Find.Execute("float")
$r.Italic = 1
$r.Collapse(0)
$r.Find.Execute("integer")
$r.Italic = 1

# ---------------------------------------------------------------------------
# 2. "External vocabularies" paragraph
# ---------------------------------------------------------------------------

$p15 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "For the XML serializations, we chose to incorporate*") {
        $p15 = $cand
        break
    }
}

$old = "For the XML serializations, we chose to incorporate two standard external vocabularies, Dublin Core and Friend of a Friend. First we set an imaginary default namespace in each XML document (" + [char]8220 + "https://example.com/maryse_and_joe_project/" + [char]8221 + ") with which to define all the elements that aren" + [char]8217 + "t given namespaces. "
$new = "For the XML serializations, we chose to incorporate bits and pieces of four authoritative external vocabularies. "

$r = $p15.Range
$r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# ---------------------------------------------------------------------------
# 3. Header name fix-up (merge the two runs / drop the spell-check proofing mark; no visible
#    text change).
# ---------------------------------------------------------------------------

foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Range.Text -like "Joseph Muller*") {
        $hr = $hdr.Range
        $hr.Find.Execute(
            "Joseph Muller & Maryse Lundering-Timpano",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "Joseph Muller & Maryse Lundering-Timpano",
            2)
    }
}

Write-Output "done"
